$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$zhTimestamp = "2016-03-09 09:42:40"
$deTimestamp = "2016-03-09 09:42:44"

$rows = @(7, 10, 11, 12, 13, 14, 15)

foreach ($r in $rows) {
    $zh.Cells.Item($r, 4).Value = $zhTimestamp
    $de.Cells.Item($r, 4).Value = $deTimestamp
}
